$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B11").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Interior.ThemeColor = 7
$ws.Range("D5").Interior.TintAndShade = 0.39997558519241921
$ws.Range("D5").HorizontalAlignment = -4108
$ws.Range("D5").VerticalAlignment = -4108
